$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 119
$ws.Range("I55").Value = 133.33333
$ws.Range("J55").Value = 97.5
$ws.Range("K55").Value = 133.33333
$ws.Range("L55").Value = 97.5
$ws.Range("M55").Value = 80.66667000000001
$ws.Range("N55").Value = -525.5

$ws.Range("H98").Value = 4360.6924
$ws.Range("I98").Value = 4244.4546
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 4244.4546
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -2746.4546
$ws.Range("N98").Value = -7996

$ws.Range("H121").Value = 704.3333
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 704.3333
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2112.9999
$ws.Range("N121").Value = -5606.9999

$ws.Range("H122").Value = 4360.6924
$ws.Range("I122").Value = 4244.4546
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 12733.3638
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -10283.3638
$ws.Range("N122").Value = -19900

$ws.Range("H137").Value = 41668756
$ws.Range("I137").Value = 1395.8334
$ws.Range("J137").Value = 83336110
$ws.Range("K137").Value = 4187.5002
$ws.Range("L137").Value = 250008330
$ws.Range("M137").Value = -1637.5002
$ws.Range("N137").Value = -250013430

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1962668.5
$ws.Range("I2").Value = 2282
$ws.Range("J2").Value = 3678006.8
$ws.Range("K2").Value = 2282
$ws.Range("L2").Value = 3678006.8
$ws.Range("M2").Value = -2169
$ws.Range("N2").Value = -3678232.8

$ws.Range("H32").Value = 4449.25
$ws.Range("I32").Value = 5031.9
$ws.Range("J32").Value = 2830.7778
$ws.Range("K32").Value = 5031.9
$ws.Range("L32").Value = 2830.7778
$ws.Range("M32").Value = -4744.9
$ws.Range("N32").Value = -3404.7778

$ws.Range("H116").Value = 1962668.5
$ws.Range("I116").Value = 2282
$ws.Range("J116").Value = 3678006.8
$ws.Range("K116").Value = 2282
$ws.Range("L116").Value = 3678006.8
$ws.Range("M116").Value = 12
$ws.Range("N116").Value = -3682594.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1962668.5
$ws.Range("I3").Value = 2282
$ws.Range("J3").Value = 3678006.8
$ws.Range("K3").Value = 2282
$ws.Range("L3").Value = 3678006.8
$ws.Range("M3").Value = -2168
$ws.Range("N3").Value = -3678234.8

$ws.Range("H105").Value = 1181.1111
$ws.Range("I105").Value = 746
$ws.Range("J105").Value = 1725
$ws.Range("K105").Value = 746
$ws.Range("L105").Value = 1725
$ws.Range("M105").Value = 1001
$ws.Range("N105").Value = -5219

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2308
$ws.Range("I31").Value = 1255.8182
$ws.Range("J31").Value = 3526.3157
$ws.Range("K31").Value = 1255.8182
$ws.Range("L31").Value = 3526.3157
$ws.Range("M31").Value = -960.8181999999999
$ws.Range("N31").Value = -4116.3157

$ws.Range("H34").Value = 2308
$ws.Range("I34").Value = 1255.8182
$ws.Range("J34").Value = 3526.3157
$ws.Range("K34").Value = 1255.8182
$ws.Range("L34").Value = 3526.3157
$ws.Range("M34").Value = -1053.8182
$ws.Range("N34").Value = -3930.3157

$ws.Range("H105").Value = 1124.3334
$ws.Range("I105").Value = 1168.625
$ws.Range("J105").Value = 1035.75
$ws.Range("K105").Value = 1168.625
$ws.Range("L105").Value = 1035.75
$ws.Range("M105").Value = 578.375
$ws.Range("N105").Value = -4529.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 475.25
$ws.Range("I10").Value = 475.25
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1425.75
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1286.75
$ws.Range("N10").ClearContents()

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

$ws.Range("H94").Value = 2561.3845
$ws.Range("I94").Value = 1932.6666
$ws.Range("J94").Value = 2750
$ws.Range("K94").Value = 5797.9998
$ws.Range("L94").Value = 8250
$ws.Range("M94").Value = -5121.9998
$ws.Range("N94").Value = -9602

$ws.Range("H95").Value = 3000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 3000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 9000
$ws.Range("N95").Value = -13118

$ws.Range("H96").Value = 35354950
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 35354950
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 106064850
$ws.Range("N96").Value = -106068968

$ws.Range("H101").Value = 6169.6
$ws.Range("I101").Value = 5000
$ws.Range("J101").Value = 8898.666999999999
$ws.Range("K101").Value = 15000
$ws.Range("L101").Value = 26696.001
$ws.Range("M101").Value = -12566
$ws.Range("N101").Value = -31564.001

$ws.Range("H105").Value = 908000000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 908000000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2724000000
$ws.Range("N105").Value = -2724005242

$ws.Range("H109").Value = 1622.8889
$ws.Range("I109").Value = 900.8570999999999
$ws.Range("J109").Value = 4150
$ws.Range("K109").Value = 2702.5713
$ws.Range("L109").Value = 12450
$ws.Range("M109").Value = -1662.5713
$ws.Range("N109").Value = -14530

$ws.Range("H121").Value = 35715236
$ws.Range("I121").Value = 336.4
$ws.Range("J121").Value = 55556850
$ws.Range("K121").Value = 1009.2
$ws.Range("L121").Value = 166670550
$ws.Range("M121").Value = 300.8000000000001
$ws.Range("N121").Value = -166673170

$ws.Range("H131").Value = 1615621
$ws.Range("I131").Value = 4894.385
$ws.Range("J131").Value = 2042956.6
$ws.Range("K131").Value = 14683.155
$ws.Range("L131").Value = 6128869.800000001
$ws.Range("M131").Value = -9643.155000000001
$ws.Range("N131").Value = -6138949.800000001

$ws.Range("H133").Value = 6689.684
$ws.Range("I133").Value = 4980
$ws.Range("J133").Value = 6890.8237
$ws.Range("K133").Value = 14940
$ws.Range("L133").Value = 20672.4711
$ws.Range("M133").Value = -9880
$ws.Range("N133").Value = -30792.4711

$ws.Range("H137").Value = 62850.176
$ws.Range("I137").Value = 2683.6365
$ws.Range("J137").Value = 173155.5
$ws.Range("K137").Value = 8050.9095
$ws.Range("L137").Value = 519466.5
$ws.Range("M137").Value = -2950.9095
$ws.Range("N137").Value = -529666.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1548.7407
$ws.Range("I16").Value = 1548.7407
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1548.7407
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1378.7407
$ws.Range("N16").ClearContents()

$ws.Range("H74").Value = 33742.855
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 33742.855
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 33742.855
$ws.Range("N74").Value = -35738.855

$ws.Range("H77").Value = 33742.855
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 33742.855
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 101228.565
$ws.Range("N77").Value = -111212.565

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 43215
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 43215
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 43215
$ws.Range("N123").Value = -53015

$ws.Range("H136").Value = 753.9729599999999
$ws.Range("I136").Value = 668.4074000000001
$ws.Range("J136").Value = 985
$ws.Range("K136").Value = 2005.2222
$ws.Range("L136").Value = 2955
$ws.Range("M136").Value = 544.7777999999998
$ws.Range("N136").Value = -8055
